$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
  "ärgern",
  "brauchen",
  "helfen",
  "leihen",
  "spüren",
  "sprengen",
  "streichen",
  "flehen",
  "siegen",
  "schwimmen",
  "sinken",
  "filmen",
  "zögern",
  "hoffen",
  "machen",
  "folgen",
  "deuten",
  "fangen",
  "sichern",
  "kümmern",
  "fließen",
  "stärken",
  "drehen",
  "flüstern",
  "schenken",
  "pressen",
  "wohnen",
  "stammen",
  "fallen",
  "boxen",
  "greifen",
  "tropfen",
  "lösen",
  "altern",
  "bergen",
  "passen",
  "merken",
  "schrecken",
  "bremsen",
  "ruhen",
  "leisten",
  "werfen",
  "fahren",
  "erben",
  "suchen",
  "orten",
  "wüten",
  "zeigen",
  "spielen",
  "fragen",
  "ehren",
  "formen",
  "weichen",
  "lassen",
  "liefern",
  "zahlen",
  "rufen",
  "mauern",
  "schaden",
  "graben",
  "meistern",
  "kehren",
  "biegen",
  "ändern",
  "wehtun",
  "wecken",
  "klettern",
  "sorgen",
  "proben",
  "holen",
  "pflanzen",
  "garen",
  "bellen",
  "wundern",
  "wählen",
  "beten",
  "meinen",
  "trauen",
  "scheitern",
  "achten",
  "schließen",
  "tollen",
  "quälen",
  "morden",
  "stehlen",
  "arten",
  "trösten",
  "wachsen",
  "stecken",
  "enden",
  "knarren",
  "kriegen",
  "sperren",
  "spenden",
  "malen",
  "schwächen",
  "flüchten",
  "loben",
  "betteln",
  "bluten",
  "sitzen",
  "bitten",
  "lockern",
  "ächzen",
  "decken",
  "schulden",
  "heulen",
  "weinen",
  "heben",
  "räumen",
  "albern",
  "planen",
  "töten",
  "saufen",
  "zielen",
  "handeln",
  "backen",
  "posten",
  "hören",
  "sprechen",
  "lesen",
  "scheinen",
  "gründen",
  "sterben",
  "leiden",
  "lügen",
  "grüßen",
  "kosten",
  "warnen",
  "trotzen",
  "führen",
  "platzen",
  "liegen",
  "klingen",
  "kichern",
  "stürmen",
  "fällen",
  "heilen",
  "äußern",
  "reizen",
  "füllen",
  "regeln",
  "stürzen",
  "klappen",
  "fischen",
  "runden",
  "reisen",
  "schneiden",
  "seufzen",
  "sammeln",
  "klingeln",
  "dienen",
  "dringen",
  "kranken",
  "schlucken",
  "rasen",
  "dauern",
  "zünden",
  "schätzen",
  "schwören",
  "schütteln",
  "spinnen",
  "mögen",
  "wenden",
  "schmecken",
  "münzen",
  "teilen",
  "freuen",
  "gelten",
  "geben",
  "duschen",
  "trennen",
  "strahlen",
  "jubeln",
  "segeln",
  "wirken",
  "bauen",
  "werden",
  "feiern",
  "hauen",
  "schwingen",
  "ziehen",
  "irren",
  "jagen",
  "reiten",
  "treiben",
  "narren",
  "wärmen",
  "pfeifen",
  "öffnen",
  "schreiten",
  "wandern"
)

for ($i = 0; $i -lt $words.Length; $i++) {
  $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}

